$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "value" header (B1) to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Extend column A's date formatting (style) down through row 22, matching A2's style
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)

# Clear the old B2 numeric value - it moves to B3 in the new layout
$ws.Range("B2").ClearContents()

# Write the full date series (year-end dates) into A2:A22
$dates = @(38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657, 46022)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Write the year-over-year values into B3:B21
$values = @(11.51866285751828, 6.007290114363029, 4.974642158654663, -8.854978371166311, 13.37451445936573, 8.397157288535361, 3.335027872728791, 0.1226029201931267, 6.350193621343236, 5.122443676600863, 3.289893304242164, 5.201177892156705, 3.567305512643082, 2.482498593966143, -8.834100858716409, 2.402478842946154, 0.9490257960172555, -2.275419501954867, -0.139401726460564)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
